# Update "Förändrad" (Changed) date column (C) for rows 2-13 from
# serial 45233 (2023-11-03) to serial 45243 (2023-11-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value2 = 45243
    }
}
